$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text, and many of the new values look
# like numbers (e.g. "321.49", "0.02330", "28.560.71"). Writing such a
# string into a cell still in the default "General" format makes Excel COM
# silently reinterpret it as a floating point number (dropping meaningful
# trailing zeros and changing the stored type), so every Price cell that is
# about to be rewritten is switched to Text format first, cell by cell (a
# multi-area Range(...) union does not reliably propagate NumberFormat to
# every area here).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.560.71"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.953.30"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.49"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4762"
$ws.Range("E7").Value = "  -5.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4033"
$ws.Range("E8").Value = "  -4.70%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08487"
$ws.Range("E10").Value = "  -6.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.059"
$ws.Range("E11").Value = "  -5.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.12"
$ws.Range("E12").Value = "  -5.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.960.08"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.601"
$ws.Range("E14").Value = "  -5.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.205"
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.016"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001077"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.93"
$ws.Range("E18").Value = "  -5.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06616"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  -5.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.014"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.814"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.593.83"
$ws.Range("E23").Value = "  -3.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.51"
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.163.78"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.09"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.10"
$ws.Range("E28").Value = "  -3.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.921"
$ws.Range("E29").Value = "  -7.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.156"
$ws.Range("E30").Value = "  -6.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.54"
$ws.Range("E31").Value = "  -3.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9972"
$ws.Range("E32").Value = "  -5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09569"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.671"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.439"
$ws.Range("E35").Value = "  -8.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.599"
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02330"
$ws.Range("E37").Value = "  -6.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06224"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.258"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.721"
$ws.Range("E40").Value = "  -6.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6230"
$ws.Range("E41").Value = "  -5.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("E42").Value = "  -5.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.013"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1923"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5943"
$ws.Range("E46").Value = "  -6.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.00"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.064"
$ws.Range("E48").Value = "  -6.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.411"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06820"
$ws.Range("E51").Value = "  -2.54%  "
